$d = $word.ActiveDocument

# Shared namespace declaration used for every InsertXML payload (Flat-OPC wrapper).
function New-FlatOpc([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) "Définitions" heading: swap which bookmark (_Toc38274974 / _Toc58335093)
#    wraps just the word vs. the word + trailing space.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(40)
$p1Xml = '<w:p w14:paraId="5AC594A1" w14:textId="77777777" w:rsidR="003A2865" w:rsidRDefault="007F41F2"><w:pPr><w:pStyle w:val="Titre1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:bookmarkStart w:id="5" w:name="_Toc58335093"/><w:bookmarkStart w:id="6" w:name="_Toc38274974"/><w:r><w:t>Définitions</w:t></w:r><w:bookmarkEnd w:id="5"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkEnd w:id="6"/></w:p>'
$p1.Range.InsertXML((New-FlatOpc $p1Xml))

# ---------------------------------------------------------------------------
# 2) "Prestations attendues" heading: same bookmark-order swap pattern
#    (_Toc38274975 / _Toc58335099).
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(53)
$p2Xml = '<w:p w14:paraId="7C0E6840" w14:textId="77777777" w:rsidR="003A2865" w:rsidRDefault="007F41F2"><w:pPr><w:pStyle w:val="Titre1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:bookmarkStart w:id="12" w:name="_Toc58335099"/><w:bookmarkStart w:id="13" w:name="_Toc38274975"/><w:r><w:t>Prestation</w:t></w:r><w:r w:rsidR="00C643CA"><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> attendues</w:t></w:r><w:bookmarkEnd w:id="12"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkEnd w:id="13"/></w:p>'
$p2.Range.InsertXML((New-FlatOpc $p2Xml))

# ---------------------------------------------------------------------------
# 3) Planning paragraph + the hyperlink paragraph that followed it are
#    replaced by a single paragraph of plain text (no more hyperlink).
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(59)
$p4 = $d.Paragraphs(60)
$mergedRange = $d.Range($p3.Range.Start, $p4.Range.End)
$p3Xml = '<w:p><w:r><w:t xml:space="preserve">Le planning est fait et régulièrement mis à jour sur </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ludus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-net.</w:t></w:r></w:p>'
$mergedRange.InsertXML((New-FlatOpc $p3Xml))
